$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.378228666666666
$ws.Range("H2").Value = 10.134686
$ws.Range("I2").Value = 0.1977735183221229
$ws.Range("J2").Value = 0.1977735183221229
$ws.Range("M2").Value = 0.1375076666666667
$ws.Range("N2").Value = 0.412523
$ws.Range("O2").Value = 0.0640991855118124
$ws.Range("P2").Value = 0.0640991855118124
$ws.Range("Q2").Value = 0.4645323414197777
$ws.Range("R2").Value = 4.180791072778
$ws.Range("S2").Value = 0.01267712144025359
$ws.Range("T2").Value = 0.01267712144025359

$ws.Range("G3").Value = 3.378228666666666
$ws.Range("H3").Value = 10.134686
$ws.Range("I3").Value = 0.1977735183221229
$ws.Range("J3").Value = 0.1977735183221229
$ws.Range("O3").Value = 0.8356360102664855
$ws.Range("P3").Value = 0.8356360102664855
$ws.Range("Q3").Value = 6.05592644156511
$ws.Range("R3").Value = 54.50333797408599
$ws.Range("S3").Value = 0.1652666737870645
$ws.Range("T3").Value = 0.1652666737870645

$ws.Range("G4").Value = 3.378228666666666
$ws.Range("H4").Value = 10.134686
$ws.Range("I4").Value = 0.1977735183221229
$ws.Range("J4").Value = 0.1977735183221229
$ws.Range("M4").Value = 0.2150913333333333
$ws.Range("N4").Value = 0.645274
$ws.Range("O4").Value = 0.1002648042217022
$ws.Range("P4").Value = 0.1002648042217021
$ws.Range("Q4").Value = 0.7266277082182221
$ws.Range("R4").Value = 6.539649373963999
$ws.Range("S4").Value = 0.01982972309480488
$ws.Range("T4").Value = 0.01982972309480488

$ws.Range("I5").Value = 0.6780480282745078
$ws.Range("J5").Value = 0.6780480282745078
$ws.Range("M5").Value = 0.1375076666666667
$ws.Range("N5").Value = 0.412523
$ws.Range("O5").Value = 0.0640991855118124
$ws.Range("P5").Value = 0.0640991855118124
$ws.Range("Q5").Value = 1.592605728216889
$ws.Range("R5").Value = 14.333451553952
$ws.Range("S5").Value = 0.0434623263502863
$ws.Range("T5").Value = 0.0434623263502863

$ws.Range("I6").Value = 0.6780480282745078
$ws.Range("J6").Value = 0.6780480282745078
$ws.Range("O6").Value = 0.8356360102664855
$ws.Range("P6").Value = 0.8356360102664855
$ws.Range("S6").Value = 0.5666013491163668
$ws.Range("T6").Value = 0.5666013491163668

$ws.Range("I7").Value = 0.6780480282745078
$ws.Range("J7").Value = 0.6780480282745078
$ws.Range("M7").Value = 0.2150913333333333
$ws.Range("N7").Value = 0.645274
$ws.Range("O7").Value = 0.1002648042217022
$ws.Range("P7").Value = 0.1002648042217021
$ws.Range("Q7").Value = 2.491175203975111
$ws.Range("R7").Value = 22.420576835776
$ws.Range("S7").Value = 0.0679843528078547
$ws.Range("T7").Value = 0.06798435280785468

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.121129333333334
$ws.Range("H8").Value = 6.363388
$ws.Range("I8").Value = 0.1241784534033691
$ws.Range("J8").Value = 0.1241784534033691
$ws.Range("M8").Value = 0.1375076666666667
$ws.Range("N8").Value = 0.412523
$ws.Range("O8").Value = 0.0640991855118124
$ws.Range("P8").Value = 0.0640991855118124
$ws.Range("Q8").Value = 0.2916715453248889
$ws.Range("R8").Value = 2.625043907924
$ws.Range("S8").Value = 0.00795973772127251
$ws.Range("T8").Value = 0.00795973772127251

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.121129333333334
$ws.Range("H9").Value = 6.363388
$ws.Range("I9").Value = 0.1241784534033691
$ws.Range("J9").Value = 0.1241784534033691
$ws.Range("O9").Value = 0.8356360102664855
$ws.Range("P9").Value = 0.8356360102664855
$ws.Range("Q9").Value = 3.802407854287555
$ws.Range("R9").Value = 34.221670688588
$ws.Range("S9").Value = 0.1037679873630541
$ws.Range("T9").Value = 0.1037679873630541

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.121129333333334
$ws.Range("H10").Value = 6.363388
$ws.Range("I10").Value = 0.1241784534033691
$ws.Range("J10").Value = 0.1241784534033691
$ws.Range("M10").Value = 0.2150913333333333
$ws.Range("N10").Value = 0.645274
$ws.Range("O10").Value = 0.1002648042217022
$ws.Range("P10").Value = 0.1002648042217021
$ws.Range("Q10").Value = 0.4562365364791112
$ws.Range("R10").Value = 4.106128828312
$ws.Range("S10").Value = 0.01245072831904257
$ws.Range("T10").Value = 0.01245072831904257
